$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set/modify cell values
$ws.Range("A52").Value = 111901584
$ws.Range("B52").Value = 56430
$ws.Range("Q52").Value = 478211
$ws.Range("R52").Value = 7035067
$ws.Range("AC52").Value = "ringhack"
$ws.Range("A53").Value = 111901585
$ws.Range("B53").Value = 56430
$ws.Range("E53").Value = 100109
$ws.Range("F53").Value = "Tretåig hackspett"
$ws.Range("G53").Value = "Picoides tridactylus"
$ws.Range("H53").Value = "(Linnaeus, 1758)"
$ws.Range("K53").Value = ""
$ws.Range("L53").Value = ""
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = ""
$ws.Range("Q53").Value = 478339
$ws.Range("R53").Value = 7035076
$ws.Range("AC53").Value = "ringhack äldre"
$ws.Range("A54").Value = 111901619
$ws.Range("B54").Value = 85183
$ws.Range("E54").Value = 249278
$ws.Range("F54").Value = "Barrviolspindling"
$ws.Range("G54").Value = "Cortinarius harcynicus"
$ws.Range("H54").Value = "(Pers.) M.M.Moser"
$ws.Range("Q54").Value = 478523
$ws.Range("R54").Value = 7034651
$ws.Range("A55").Value = 112102606
$ws.Range("B55").Value = 89033
$ws.Range("E55").Value = 3286
$ws.Range("F55").Value = "Flattoppad klubbsvamp"
$ws.Range("G55").Value = "Clavariadelphus truncatus"
$ws.Range("H55").Value = "(Quél.) Donk"
$ws.Range("P55").Value = "Storbäcken (Storbäcken), Jmt"
$ws.Range("Q55").Value = 478088
$ws.Range("R55").Value = 7035319
$ws.Range("S55").Value = 5
$ws.Range("Y55").Value = "'2023-09-15"
$ws.Range("Z55").Value = "'10:16"
$ws.Range("AA55").Value = "'2023-09-15"
$ws.Range("AB55").Value = "'10:16"
$ws.Range("AW55").Value = "Jonny Daborg"
$ws.Range("AX55").Value = "Jonny Daborg"
$ws.Range("A56").Value = 111901546
$ws.Range("B56").Value = 56430
$ws.Range("Q56").Value = 477668
$ws.Range("R56").Value = 7033374
$ws.Range("AC56").Value = "ringhack äldre"
$ws.Range("A57").Value = 111901587
$ws.Range("B57").Value = 56575
$ws.Range("E57").Value = 103021
$ws.Range("F57").Value = "Talltita"
$ws.Range("G57").Value = "Poecile montanus"
$ws.Range("H57").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I57").Value = "'2"
$ws.Range("N57").Value = "observerad"
$ws.Range("Q57").Value = 477611
$ws.Range("R57").Value = 7033311
$ws.Range("A58").Value = 111901519
$ws.Range("B58").Value = 86357
$ws.Range("E58").Value = 4412
$ws.Range("F58").Value = "Äggvaxskivling"
$ws.Range("G58").Value = "Hygrophorus karstenii"
$ws.Range("H58").Value = "Sacc. & Cub."
$ws.Range("Q58").Value = 477765
$ws.Range("R58").Value = 7033404
$ws.Range("A59").Value = 111901548
$ws.Range("B59").Value = 56430
$ws.Range("Q59").Value = 477476
$ws.Range("R59").Value = 7033385
$ws.Range("A60").Value = 111901547
$ws.Range("B60").Value = 56430
$ws.Range("E60").Value = 100109
$ws.Range("F60").Value = "Tretåig hackspett"
$ws.Range("G60").Value = "Picoides tridactylus"
$ws.Range("H60").Value = "(Linnaeus, 1758)"
$ws.Range("K60").Value = ""
$ws.Range("L60").Value = ""
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = ""
$ws.Range("Q60").Value = 477524
$ws.Range("R60").Value = 7033330
$ws.Range("AC60").Value = "ringhack"
$ws.Range("B61").Value = 56430
$ws.Range("A62").Value = 111901518
$ws.Range("B62").Value = 86357
$ws.Range("E62").Value = 4412
$ws.Range("F62").Value = "Äggvaxskivling"
$ws.Range("G62").Value = "Hygrophorus karstenii"
$ws.Range("H62").Value = "Sacc. & Cub."
$ws.Range("Q62").Value = 477674
$ws.Range("R62").Value = 7033500
$ws.Range("A63").Value = 111901549
$ws.Range("B63").Value = 56430
$ws.Range("Q63").Value = 477464
$ws.Range("R63").Value = 7033364
$ws.Range("AC63").Value = "ringhack färska"
$ws.Range("A64").Value = 111901551
$ws.Range("B64").Value = 56430
$ws.Range("Q64").Value = 477433
$ws.Range("R64").Value = 7033429
$ws.Range("AC64").Value = "ringhack"
$ws.Range("A65").Value = 111901550
$ws.Range("B65").Value = 56430
$ws.Range("E65").Value = 100109
$ws.Range("F65").Value = "Tretåig hackspett"
$ws.Range("G65").Value = "Picoides tridactylus"
$ws.Range("H65").Value = "(Linnaeus, 1758)"
$ws.Range("K65").Value = ""
$ws.Range("L65").Value = ""
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = ""
$ws.Range("Q65").Value = 477473
$ws.Range("R65").Value = 7033404
$ws.Range("AC65").Value = "ringhack äldre"
$ws.Range("A66").Value = 111901618
$ws.Range("B66").Value = 85183
$ws.Range("E66").Value = 249278
$ws.Range("F66").Value = "Barrviolspindling"
$ws.Range("G66").Value = "Cortinarius harcynicus"
$ws.Range("H66").Value = "(Pers.) M.M.Moser"
$ws.Range("I66").Value = ""
$ws.Range("Q66").Value = 477471
$ws.Range("R66").Value = 7033412
$ws.Range("A67").Value = 111901544
$ws.Range("B67").Value = 56430
$ws.Range("E67").Value = 100109
$ws.Range("F67").Value = "Tretåig hackspett"
$ws.Range("G67").Value = "Picoides tridactylus"
$ws.Range("H67").Value = "(Linnaeus, 1758)"
$ws.Range("L67").Value = ""
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = ""
$ws.Range("P67").Value = "Landvågen, Jmt"
$ws.Range("Q67").Value = 477639
$ws.Range("R67").Value = 7033515
$ws.Range("S67").Value = 10
$ws.Range("Y67").Value = "'2023-09-04"
$ws.Range("AA67").Value = "'2023-09-04"
$ws.Range("AC67").Value = "ringhack äldre"
$ws.Range("AW67").Value = "Benny Öwre"
$ws.Range("AX67").Value = "Benny Öwre"

# Clear cells that are removed in the target state
$ws.Range("L55").ClearContents()
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("AC55").ClearContents()
$ws.Range("AC57").ClearContents()
$ws.Range("K58").ClearContents()
$ws.Range("L58").ClearContents()
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("AC58").ClearContents()
$ws.Range("K62").ClearContents()
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("AC62").ClearContents()
$ws.Range("K66").ClearContents()
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("Z67").ClearContents()
$ws.Range("AB67").ClearContents()
